$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 57-64 (participant "26" rows), which shifts subsequent rows up.
$ws.Range("A57:L64").EntireRow.Delete()
